$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 34, shifting existing rows 34-91 down to 36-93
$ws.Range("A34:A35").EntireRow.Insert()

# New row 34: Kiwi Hayward Primera, 2021-08-05 (44413)
$ws.Cells.Item(34,1).Value = 7
$ws.Cells.Item(34,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(34,3).Value = "Ñuble"
$ws.Cells.Item(34,4).Value = 44413
$ws.Cells.Item(34,5).Value = 16
$ws.Cells.Item(34,6).Value = "Fruta"
$ws.Cells.Item(34,7).Value = 100101
$ws.Cells.Item(34,8).Value = "Berries"
$ws.Cells.Item(34,9).Value = 100101007
$ws.Cells.Item(34,10).Value = "Kiwi"
$ws.Cells.Item(34,11).Value = "Hayward"
$ws.Cells.Item(34,12).Value = "Primera"
$ws.Cells.Item(34,13).Value = 200
$ws.Cells.Item(34,14).Value = 12500
$ws.Cells.Item(34,15).Value = 13000
$ws.Cells.Item(34,16).Value = 12750
$ws.Cells.Item(34,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(34,18).Value = "Provincia de Curicó"
$ws.Cells.Item(34,19).Value = 708
$ws.Cells.Item(34,20).Value = 18

# New row 35: Kiwi Hayward Segunda, 2021-08-05 (44413)
$ws.Cells.Item(35,1).Value = 7
$ws.Cells.Item(35,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35,3).Value = "Ñuble"
$ws.Cells.Item(35,4).Value = 44413
$ws.Cells.Item(35,5).Value = 16
$ws.Cells.Item(35,6).Value = "Fruta"
$ws.Cells.Item(35,7).Value = 100101
$ws.Cells.Item(35,8).Value = "Berries"
$ws.Cells.Item(35,9).Value = 100101007
$ws.Cells.Item(35,10).Value = "Kiwi"
$ws.Cells.Item(35,11).Value = "Hayward"
$ws.Cells.Item(35,12).Value = "Segunda"
$ws.Cells.Item(35,13).Value = 80
$ws.Cells.Item(35,14).Value = 11000
$ws.Cells.Item(35,15).Value = 11000
$ws.Cells.Item(35,16).Value = 11000
$ws.Cells.Item(35,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(35,18).Value = "Provincia de Curicó"
$ws.Cells.Item(35,19).Value = 611
$ws.Cells.Item(35,20).Value = 18
